$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.561.98'
$ws.Range("E2").Value = '  +1.16%  '
$ws.Range("D3").Value = '1.854.87'
$ws.Range("E3").Value = '  +0.35%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9997'
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '233.65'
$ws.Range("E5").Value = '  +0.39%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.0000'
$ws.Range("E6").Value = '  -0.06%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4732'
$ws.Range("E7").Value = '  +0.63%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2749'
$ws.Range("E8").Value = '  +1.40%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06319'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '17.77'
$ws.Range("E10").Value = '  +9.18%  '
$ws.Range("D11").Value = '1.876.15'
$ws.Range("E11").Value = '  +1.42%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07446'
$ws.Range("E12").Value = '  +0.25%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.999'
$ws.Range("E13").Value = '  +1.31%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '84.64'
$ws.Range("E14").Value = '  -0.48%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6264'
$ws.Range("E15").Value = '  +0.04%  '
$ws.Range("D16").Value = '30.527.79'
$ws.Range("E16").Value = '  +1.20%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '246.21'
$ws.Range("E17").Value = '  +7.74%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9996'
$ws.Range("E18").Value = '  -0.03%  '
$ws.Range("E19").Value = '  +1.06%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007335'
$ws.Range("E20").Value = '  +0.08%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9985'
$ws.Range("E21").Value = '  -0.16%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.946'
$ws.Range("E22").Value = '  +0.40%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.945'
$ws.Range("E23").Value = '  +0.48%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.143'
$ws.Range("E24").Value = '  -0.89%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '162.86'
$ws.Range("E25").Value = '  -2.21%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '18.01'
$ws.Range("E26").Value = '  +1.00%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.878'
$ws.Range("E27").Value = '  +0.37%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.1020'
$ws.Range("E28").Value = '  -0.86%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.360'
$ws.Range("E29").Value = '  -1.29%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.009'
$ws.Range("E30").Value = '  -2.56%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.836'
$ws.Range("E31").Value = '  -0.93%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.04842'
$ws.Range("E32").Value = '  -0.89%  '
$ws.Range("E33").Value = '  -1.19%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7036'
$ws.Range("E34").Value = '  -0.49%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.690'
$ws.Range("E35").Value = '  -0.28%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.01897'
$ws.Range("E36").Value = '  +2.12%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.675'
$ws.Range("E37").Value = '  +1.68%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.009'
$ws.Range("E38").Value = '  +3.08%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.8761'
$ws.Range("E39").Value = '  -3.17%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '106.90'
$ws.Range("E40").Value = '  +1.92%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9998'
$ws.Range("E41").Value = '  +0.16%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.538'
$ws.Range("E42").Value = '  -0.31%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4062'
$ws.Range("E43").Value = '  -0.48%  '
$ws.Range("E44").Value = '  +1.80%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '62.68'
$ws.Range("E45").Value = '  +3.96%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1211'
$ws.Range("E46").Value = '  +1.73%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '33.64'
$ws.Range("E47").Value = '  +1.59%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.512'
$ws.Range("E48").Value = '  -1.10%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05537'
$ws.Range("E49").Value = '  -0.47%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.354'
$ws.Range("E50").Value = '  -2.17%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3687'
$ws.Range("E51").Value = '  +0.56%  '
